$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

# Add the two new rows (order matches the target shared-string insertion order:
# Path column first for both rows, then ContentType column for both rows).
$ws.Range("A3").Value = "/espanol"
$ws.Range("A4").Value = "/sites/nano"
$ws.Range("B3").Value = "Home/Landing Page"
$ws.Range("B4").Value = "Home/Landing Page"

# Widen column B to fit the new, longer "Home/Landing Page" content.
$ws.Columns.Item(2).ColumnWidth = 18.17

# Make HomePage the active/selected sheet (was LandingPage before the edit),
# and leave the selection on A2.
$ws.Activate()
$ws.Range("A2").Select()
